# Consent Decision Reporting mapping sheet: collapse the "Consent" code
# list down to a single combined "Codes: ..." note, and drop the
# Booking Facility Name / Booking Facility Number rows entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Booking Facility Name" / "Booking Facility Number" rows
# (originally rows 17-18, right after "Booking Number").
$ws.Rows("17:18").Delete() | Out-Null

# The old "Consent Denied" / "Inmate Never Seen" / "Consent Not Obtained"
# rows (originally 24-26, now 22-24 after the delete above) only carried a
# column-B code label with no A/C text - fold them away too.
$ws.Rows("22:24").Delete() | Out-Null

# The remaining "Consent Decision Code" row (now row 21) keeps its single
# "Consent Granted" label in B - replace it with the combined code note.
$ws.Cells.Item(21, 2).Value2 = "Codes: Consent Granted; Consent Denied"

# Restore the sheet's last selection so the saved view matches.
$ws.Range("B29").Select() | Out-Null
